# Sync attendance_reports, modules_schedules, and assets from main repo - 2026-01-04 10:15:24
#
# The upstream "Recorded By" column (G) listed the recording users in an
# inconsistent order. This normalizes the ordering for the specific
# combinations that changed upstream by swapping the first two names in the
# comma-separated list, leaving every other combination (and every other
# column) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Exact-match replacements for the "Recorded By" values that changed.
$replacements = @{
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
    "backup@backdoor.com, system, System" = "system, backup@backdoor.com, System"
}

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G = "Recorded By"
    $current = $cell.Value

    if ($null -ne $current -and $replacements.ContainsKey($current)) {
        $cell.Value = $replacements[$current]
    }
}
